# Enhance navigation menu accessibility
# Adds a new data row (row 72) to each of the four worksheets, mirroring
# the existing row layout (columns A-I).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = 1
        A = [double]"45760.39731667824"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 13
    },
    @{
        Sheet = 2
        A = [double]"45760.25136153935"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 14
    },
    @{
        Sheet = 3
        A = [double]"45760.39202359954"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 3
    },
    @{
        Sheet = 4
        A = [double]"45760.45250971065"
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 358
        I = 3
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)

    $ws.Range("A72").Value = $rowData.A
    $ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B72").Value = $rowData.B
    $ws.Range("C72").Value = $rowData.C
    $ws.Range("D72").Value = $rowData.D
    $ws.Range("E72").Value = $rowData.E
    $ws.Range("F72").Value = $rowData.F
    $ws.Range("G72").Value = $rowData.G
    $ws.Range("H72").Value = $rowData.H
    $ws.Range("I72").Value = $rowData.I
}
